$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kaif Khan")

$ws.Range("A2").Value = "The Book Shop"
$ws.Range("A3").Value = "Like a Love Song"

$ws.Activate()
$ws.Range("A3").Select()
